# Updated cryptos list on Sun Jan 14 17:10:22 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to be treated as text so numeric-looking strings
    # (e.g. "306.04") are not coerced into floating point numbers,
    # then restore the default "Normal" style so no stray number
    # format / quote-prefix style is left behind.
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "42.896.95"
Set-TextValue "E2" "  +0.08%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.534.66"
Set-TextValue "E3" "  -0.83%  "

# Row 4 - TetherUSD
Set-TextValue "D4" "0.999"
Set-TextValue "E4" "  -0.07%  "

# Row 5 - BNB
Set-TextValue "D5" "306.04"
Set-TextValue "E5" "  +1.60%  "

# Row 6 - Solana
Set-TextValue "D6" "101.15"
Set-TextValue "E6" "  +7.50%  "

# Row 7 - XRP
Set-TextValue "D7" "0.579"
Set-TextValue "E7" "  +1.18%  "

# Row 8 - USDC
Set-TextValue "E8" "  +0.00%  "

# Row 9 - Cardano
Set-TextValue "D9" "0.547"
Set-TextValue "E9" "  +0.70%  "

# Row 10 - Avalanche
Set-TextValue "D10" "37.56"
Set-TextValue "E10" "  +3.48%  "

# Row 11 - Dogecoin
Set-TextValue "D11" "0.0818"
Set-TextValue "E11" "  +0.99%  "

# Row 12 - Polkadot
Set-TextValue "D12" "7.61"
Set-TextValue "E12" "  -1.52%  "

# Row 13 - TRON
Set-TextValue "E13" "  -0.52%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextValue "D14" "2.920.87"
Set-TextValue "E14" "  -1.01%  "

# Row 15 & 16 - swap WrappedEther and Chainlink entries with new values
Set-TextValue "B15" "Chainlink"
Set-TextValue "C15" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D15" "15.22"
Set-TextValue "E15" "  +7.54%  "

Set-TextValue "B16" "WrappedEther"
Set-TextValue "C16" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D16" "2.477.21"
Set-TextValue "E16" "  -5.09%  "

# Row 17 - Polygon
Set-TextValue "D17" "0.872"
Set-TextValue "E17" "  -0.44%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "42.927.05"
Set-TextValue "E18" "  +0.05%  "

# Row 19 - InternetComputer(DFINITY)
Set-TextValue "D19" "13.24"
Set-TextValue "E19" "  +4.49%  "

# Row 20 - ShibaInu
Set-TextValue "D20" "0.0₃0988"
Set-TextValue "E20" "  +0.25%  "

# Row 21 - Uniswap
Set-TextValue "E21" "  -0.39%  "

# Row 22 - Litecoin
Set-TextValue "D22" "71.64"
Set-TextValue "E22" "  +0.36%  "

# Row 23 - BitcoinCash
Set-TextValue "D23" "253.90"
Set-TextValue "E23" "  +0.67%  "

# Row 24 - Pancake...
Set-TextValue "D24" "2.94"
Set-TextValue "E24" "  +0.18%  "

# Row 25
Set-TextValue "E25" "  -2.60%  "

# Row 26
Set-TextValue "D26" "27.28"
Set-TextValue "E26" "  -4.81%  "

# Row 27
Set-TextValue "E27" "  +0.19%  "

# Row 28 - Toncoin
Set-TextValue "D28" "2.33"
Set-TextValue "E28" "  +9.07%  "

# Row 29 - Cosmos
Set-TextValue "D29" "10.37"
Set-TextValue "E29" "  +1.73%  "

# Row 30 - InjectiveProtocol
Set-TextValue "D30" "38.75"
Set-TextValue "E30" "  +5.41%  "

# Row 31 - Filecoin
Set-TextValue "D31" "6.16"
Set-TextValue "E31" "  +1.75%  "

# Row 32 - Monero
Set-TextValue "E32" "  +2.32%  "

# Row 33 - ARBITRUM
Set-TextValue "E33" "  -0.79%  "

# Row 34 - Hedera
Set-TextValue "D34" "0.0797"
Set-TextValue "E34" "  +0.15%  "

# Row 35 - LidoDAOToken
Set-TextValue "D35" "3.29"
Set-TextValue "E35" "  -2.20%  "

# Row 36 - WEMIXToken
Set-TextValue "E36" "  -3.46%  "

# Row 37 - Celestia
Set-TextValue "D37" "18.36"
Set-TextValue "E37" "  +3.11%  "

# Row 39 - Stellar
Set-TextValue "E39" "  +0.13%  "

# Row 40 - EnergySwap
Set-TextValue "D40" "23.77"
Set-TextValue "E40" "  +2.77%  "

# Row 41 - NEARProtocol
Set-TextValue "E41" "  +2.93%  "

# Row 42 - ApeXProtocol
Set-TextValue "E42" "  -0.38%  "

# Row 43 - RenderToken
Set-TextValue "E43" "  +0.74%  "

# Row 44 - VeChain
Set-TextValue "D44" "0.0305"
Set-TextValue "E44" "  -1.30%  "

# Row 45 - FirstDigitalUSD
Set-TextValue "D45" "0.997"
Set-TextValue "E45" "  -0.17%  "

# Row 46 - Maker
Set-TextValue "D46" "2.045.45"
Set-TextValue "E46" "  -2.64%  "

# Row 47 - BitcoinSV
Set-TextValue "D47" "86.28"
Set-TextValue "E47" "  +1.94%  "

# Row 48 - FraxShare
Set-TextValue "D48" "8.98"
Set-TextValue "E48" "  -2.75%  "

# Row 49 - RocketPoolETH
Set-TextValue "D49" "2.779.53"
Set-TextValue "E49" "  -0.95%  "

# Row 50 - Algorand
Set-TextValue "E50" "  +2.08%  "

# Row 51 - Aave
Set-TextValue "D51" "103.18"
Set-TextValue "E51" "  -2.29%  "
